# Update latest output (run 35)

$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 592.3275247500002
$schedule.Range("F2").Value = 13.05836694775133
$schedule.Range("E3").Value = 357.79152825
$schedule.Range("F3").Value = 23.66346086309524

# --- Sheet: Detailed ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B19").Value = 56.97989
$detailed.Range("B20").Value = 47.31837

$detailed.Range("B21").Value = -11.03101
$detailed.Range("C21").Value = "historical"

$detailed.Range("B22").Value = -8.94445
$detailed.Range("C22").Value = "historical"

$detailed.Range("B23").Value = 36.06

$detailed.Range("B25").Value = -16.24859
$detailed.Range("B26").Value = -16.86993
$detailed.Range("B27").Value = -17.10346
$detailed.Range("B28").Value = -15.51447
$detailed.Range("B29").Value = -17.0409
$detailed.Range("B30").Value = -22.18182
$detailed.Range("B31").Value = -17.11931
$detailed.Range("B32").Value = -16.88235
$detailed.Range("B33").Value = -16.88892
$detailed.Range("B34").Value = 18.11384
$detailed.Range("B35").Value = 10.27661
$detailed.Range("B36").Value = -10.86954
$detailed.Range("B37").Value = -9.41533
$detailed.Range("B38").Value = -9.25382
$detailed.Range("B39").Value = -3.05055
$detailed.Range("B40").Value = -0.5508

$detailed.Range("B42").Value = 29.85322
$detailed.Range("B43").Value = 29.85322
$detailed.Range("B44").Value = 22.01959
$detailed.Range("B45").Value = 57.09

$detailed.Range("B49").Value = 56.98
